$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers for U1:W1
$ws.Range("U1").Value = "GG Expenditure"
$ws.Range("V1").Value = "GG Revenue"
$ws.Range("W1").Value = "GG Balance"

# Rows 2-21: create empty (but present) cells in U:W, matching the unstyled
# placeholder cells already used throughout this sheet for "no data yet" slots.
for ($r = 2; $r -le 21; $r++) {
    $ws.Range("U$r").Style = "Normal"
    $ws.Range("V$r").Style = "Normal"
    $ws.Range("W$r").Style = "Normal"
}

# Rows 22-122: the actual GG Expenditure / GG Revenue / GG Balance data
$ggData = @{
    22 = @(7374, 8558, 1184)
    23 = @(7572, 10525, 2953)
    24 = @(8378, 9145, 767)
    25 = @(9841, 10204, 363)
    26 = @(8273, 9060, 787)
    27 = @(9269, 11209, 1940)
    28 = @(9688, 9246, -441)
    29 = @(12050, 10936, -1115)
    30 = @(10095, 9544, -551)
    31 = @(10572, 10830, 257)
    32 = @(10903, 11059, 156)
    33 = @(13145, 12577, -569)
    34 = @(10790, 10502, -287)
    35 = @(11456, 11762, 306)
    36 = @(11612, 11253, -359)
    37 = @(13695, 14539, 845)
    38 = @(11257, 12056, 799)
    39 = @(12134, 12847, 714)
    40 = @(12812, 12291, -521)
    41 = @(15084, 16121, 1037)
    42 = @(12566, 13374, 808)
    43 = @(13309, 13054, -255)
    44 = @(13943, 13999, 57)
    45 = @(16421, 18485, 2064)
    46 = @(13658, 15585, 1928)
    47 = @(14650, 14226, -424)
    48 = @(15650, 15107, -543)
    49 = @(18105, 22276, 4171)
    50 = @(15898, 17331, 1433)
    51 = @(16448, 14855, -1593)
    52 = @(17891, 15741, -2150)
    53 = @(19999, 22837, 2838)
    54 = @(18240, 16423, -1817)
    55 = @(18671, 14106, -4565)
    56 = @(19183, 14907, -4276)
    57 = @(21736, 19227, -2509)
    58 = @(18366, 13450, -4916)
    59 = @(21284, 13245, -8038)
    60 = @(20137, 13518, -6619)
    61 = @(19674, 15731, -3943)
    62 = @(29060, 12535, -16526)
    63 = @(19485, 13167, -6318)
    64 = @(26487, 12876, -13611)
    65 = @(33595, 16286, -17309)
    66 = @(17655, 13225, -4429)
    67 = @(17591, 13653, -3938)
    68 = @(25682, 14290, -11392)
    69 = @(20187, 16625, -3561)
    70 = @(18420, 13694, -4728)
    71 = @(18674, 14651, -4023)
    72 = @(18201, 14773, -3428)
    73 = @(19371, 16685, -2686)
    74 = @(18167, 14109, -4057)
    75 = @(17256, 15365, -1892)
    76 = @(17871, 14844, -3027)
    77 = @(19581, 17077, -2505)
    78 = @(17804, 15034, -2770)
    79 = @(17835, 16420, -1416)
    80 = @(18060, 16003, -2056)
    81 = @(19623, 18810, -814)
    82 = @(18328, 16075, -2253)
    83 = @(18067, 17527, -540)
    84 = @(18154, 16763, -1391)
    85 = @(21881, 20693, -1189)
    86 = @(18212, 16447, -1765)
    87 = @(18249, 18446, 197)
    88 = @(18935, 17188, -1747)
    89 = @(20528, 21741, 1213)
    90 = @(18517, 17438, -1079)
    91 = @(18940, 18828, -112)
    92 = @(19690, 17793, -1896)
    93 = @(20806, 22976, 2170)
    94 = @(19644, 17936, -1709)
    95 = @(20462, 20284, -178)
    96 = @(20526, 18932, -1594)
    97 = @(22304, 26087, 3784)
    98 = @(20917, 18929, -1988)
    99 = @(20877, 21721, 844)
    100 = @(21684, 20586, -1098)
    101 = @(23343, 27063, 3720)
    102 = @(21549, 17937, -3612)
    103 = @(26654, 20257, -6397)
    104 = @(26230, 19527, -6703)
    105 = @(27498, 25628, -1870)
    106 = @(25564, 19571, -5993)
    107 = @(26230, 23800, -2430)
    108 = @(26583, 24033, -2550)
    109 = @(27408, 32251, 4842)
    110 = @(24611, 24363, -248)
    111 = @(26022, 27689, 1667)
    112 = @(26881, 28814, 1934)
    113 = @(29809, 35142, 5332)
    114 = @(27040, 27670, 630)
    115 = @(28162, 29729, 1567)
    116 = @(29020, 28886, -134)
    117 = @(31761, 37698, 5937)
    118 = @(29596, 28932, -664)
    119 = @(30490, 33546, 3055)
    120 = @(31070, 46533, 15463)
    121 = @(34390, 39530, 5140)
    122 = @(30086, 30924, 838)
}
foreach ($r in $ggData.Keys) {
    $vals = $ggData[$r]
    $ws.Range("U$r").Value = $vals[0]
    $ws.Range("V$r").Value = $vals[1]
    $ws.Range("W$r").Value = $vals[2]
}

# Rows 123-145: trailing empty placeholder cells in U:W
for ($r = 123; $r -le 145; $r++) {
    $ws.Range("U$r").Style = "Normal"
    $ws.Range("V$r").Style = "Normal"
    $ws.Range("W$r").Style = "Normal"
}
